$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Move "Expected Result" column (old C) into new column D, rows 1-3 ----
$ws.Range("D1").Value = "Expected Result"
$ws.Range("D2").Value = "Submitted for grading"
$ws.Range("D3").Value = "The file 5_Unitintergration testing.pdf is too large. The maximum size you can upload is 1 MB."

# ---- Rename column B header, add column C header ----
$ws.Range("B1").Value = "File1"
$ws.Range("C1").Value = "File2"

# Old C2/C3 ("Expected Result" values) are no longer needed there - clear them
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Row2 (TC-03-01) / Row3 (TC-03-02): column B values are already correct (unchanged)

# ---- Insert new row 4 (TC-03-03) ----
$ws.Range("A4").EntireRow.Insert()
$ws.Range("A4").Value = "TC-03-03"
$ws.Range("B4").Value = "C:\\Software-Testing\\Project #1 description.pdf"
$ws.Range("C4").Value = "C:\\Software-Testing\\Project #2 description.pdf"
$ws.Range("D4").Value = "You are allowed to attach a maximum of 1 file(s) to this item"

# ---- Row 5 - TC-03-04 (previously row 4): move its Expected Result to D, clear C ----
$ws.Range("D5").Value = "Video file (MP4) filetype cannot be accepted."
$ws.Range("C5").ClearContents()

# ---- Insert new row 6 (TC-03-05) ----
$ws.Range("A6").EntireRow.Insert()
$ws.Range("A6").Value = "TC-03-05"
$ws.Range("B6").Value = "C:\\Software-Testing\\Project #1 description.pdf"
$ws.Range("C6").Value = "C:\\Software-Testing\\5_Unitintergration testing.pdf"
$ws.Range("D6").Value = "You are allowed to attach a maximum of 1 file(s) to this item"

# ---- Row 7 - TC-03-06 (previously row 5): move its Expected Result to D, clear C ----
$ws.Range("D7").Value = "Video file (MP4) filetype cannot be accepted."
$ws.Range("C7").ClearContents()

# ---- Insert new row 8 (TC-03-07) ----
$ws.Range("A8").EntireRow.Insert()
$ws.Range("A8").Value = "TC-03-07"
$ws.Range("B8").Value = "C:\\Software-Testing\\Project #1 description.pdf"
$ws.Range("C8").Value = "C:\\Software-Testing\\singleball.mp4"
$ws.Range("D8").Value = "You are allowed to attach a maximum of 1 file(s) to this item"

# ---- New row 9 (TC-03-08) - copy format from row 8 then set values ----
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A9").Value = "TC-03-08"
$ws.Range("B9").Value = "C:\\Software-Testing\\Project #1 description.pdf"
$ws.Range("C9").Value = "C:\\Software-Testing\\vipboard.mp4"
$ws.Range("D9").Value = "You are allowed to attach a maximum of 1 file(s) to this item"

# ---- Empty styled cells at C13, C14 (style matches column-A data style) ----
$ws.Range("A2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()

[void]$ws.Range("C17").Select()
